$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB4").Value = 12
$ws.Range("AE4").Value = 10
$ws.Range("AF4").Value = 11.25
$ws.Range("AI4").Value = 13.5
$ws.Range("AJ4").Value = 19.5
$ws.Range("R4").Value = 1.52
$ws.Range("T4").Value = 14
$ws.Range("U4").Value = 21
$ws.Range("W4").Value = 45
$ws.Range("Y4").Value = 27

$ws.Range("AE9").Value = 6.6
$ws.Range("AF9").Value = 14.5
$ws.Range("AG9").Value = 11.5
$ws.Range("AH9").Value = 45
$ws.Range("AI9").Value = 37
$ws.Range("G9").Value = 2.87
$ws.Range("I9").Value = 3.1
$ws.Range("T9").Value = 6.2
$ws.Range("U9").Value = 13
$ws.Range("V9").Value = 10.75
$ws.Range("W9").Value = 40
$ws.Range("X9").Value = 32

$ws.Range("J10").Value = 1.07
$ws.Range("K10").Value = 9

$ws.Range("AA11").Value = 5.8
$ws.Range("AC11").Value = 65
$ws.Range("AD11").Value = 500
$ws.Range("AE11").Value = 9
$ws.Range("AF11").Value = 19
$ws.Range("AG11").Value = 12
$ws.Range("AH11").Value = 55
$ws.Range("AI11").Value = 37
$ws.Range("AJ11").Value = 40
$ws.Range("G11").Value = 1.72
$ws.Range("H11").Value = 3.35
$ws.Range("I11").Value = 4.4
$ws.Range("L11").Value = 1.35
$ws.Range("M11").Value = 2.95
$ws.Range("N11").Value = 1.98
$ws.Range("O11").Value = 1.65
$ws.Range("P11").Value = 1.38
$ws.Range("Q11").Value = 2.47
$ws.Range("R11").Value = 1.99
$ws.Range("T11").Value = 5.2
$ws.Range("U11").Value = 6.4
$ws.Range("V11").Value = 7
$ws.Range("W11").Value = 10.75
$ws.Range("X11").Value = 12
$ws.Range("Y11").Value = 24
$ws.Range("Z11").Value = 8.25

$ws.Range("AA12").Value = 5.3
$ws.Range("AB12").Value = 12.5
$ws.Range("AC12").Value = 60
$ws.Range("AD12").Value = 450
$ws.Range("AE12").Value = 6.2
$ws.Range("AF12").Value = 9.25
$ws.Range("AG12").Value = 8
$ws.Range("AH12").Value = 19
$ws.Range("AI12").Value = 16.5
$ws.Range("AJ12").Value = 26
$ws.Range("G12").Value = 2.75
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 2.37
$ws.Range("L12").Value = 1.38
$ws.Range("M12").Value = 2.85
$ws.Range("N12").Value = 2.05
$ws.Range("O12").Value = 1.6
$ws.Range("P12").Value = 1.4
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.89
$ws.Range("S12").Value = 1.82
$ws.Range("T12").Value = 6.6
$ws.Range("U12").Value = 10.75
$ws.Range("V12").Value = 8.75
$ws.Range("W12").Value = 25
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 29
$ws.Range("Z12").Value = 8

$ws.Range("AA13").Value = 6.4
$ws.Range("AB13").Value = 16.5
$ws.Range("AC13").Value = 80
$ws.Range("AE13").Value = 11.25
$ws.Range("AF13").Value = 27
$ws.Range("AG13").Value = 16
$ws.Range("AH13").Value = 90
$ws.Range("AI13").Value = 55
$ws.Range("AJ13").Value = 55
$ws.Range("G13").Value = 1.5
$ws.Range("H13").Value = 3.7
$ws.Range("I13").Value = 5.9
$ws.Range("L13").Value = 1.32
$ws.Range("M13").Value = 3.1
$ws.Range("N13").Value = 1.93
$ws.Range("O13").Value = 1.7
$ws.Range("P13").Value = 1.37
$ws.Range("Q13").Value = 2.5
$ws.Range("R13").Value = 2.14
$ws.Range("S13").Value = 1.63
$ws.Range("T13").Value = 4.9
$ws.Range("U13").Value = 5.4
$ws.Range("V13").Value = 7.1
$ws.Range("W13").Value = 8.25
$ws.Range("X13").Value = 11
$ws.Range("Y13").Value = 26
$ws.Range("Z13").Value = 8.75

$ws.Range("L19").Value = 1.33
$ws.Range("M19").Value = 3.25
$ws.Range("N19").Value = 2.1
$ws.Range("O19").Value = 1.7

$ws.Range("L20").Value = 1.33
$ws.Range("M20").Value = 3.25

$ws.Range("J25").Value = 1.1
$ws.Range("K25").Value = 7

$ws.Range("L26").Value = 1.25
$ws.Range("M26").Value = 3.75

$ws.Range("AB32").Value = 12
$ws.Range("AD32").Value = 300
$ws.Range("AE32").Value = 12
$ws.Range("AF32").Value = 23
$ws.Range("AG32").Value = 12
$ws.Range("AH32").Value = 60
$ws.Range("AI32").Value = 32
$ws.Range("AJ32").Value = 32
$ws.Range("G32").Value = 1.62
$ws.Range("H32").Value = 3.65
$ws.Range("I32").Value = 4.55
$ws.Range("L32").Value = 1.23
$ws.Range("M32").Value = 3.75
$ws.Range("N32").Value = 1.7
$ws.Range("O32").Value = 1.93
$ws.Range("R32").Value = 1.75
$ws.Range("S32").Value = 1.96
$ws.Range("T32").Value = 6.4
$ws.Range("U32").Value = 6.9
$ws.Range("W32").Value = 10.5
$ws.Range("X32").Value = 10.5
$ws.Range("Y32").Value = 18.5
$ws.Range("Z32").Value = 11.5

$ws.Range("L35").Value = 1.2
$ws.Range("M35").Value = 4

$ws.Range("AA36").Value = 6.5
$ws.Range("AB36").Value = 15.5
$ws.Range("AC36").Value = 80
$ws.Range("AD36").Value = 700
$ws.Range("AE36").Value = 10
$ws.Range("AF36").Value = 19
$ws.Range("AG36").Value = 12.5
$ws.Range("AH36").Value = 55
$ws.Range("AI36").Value = 35
$ws.Range("AJ36").Value = 45
$ws.Range("G36").Value = 1.93
$ws.Range("I36").Value = 3.65
$ws.Range("L36").Value = 1.32
$ws.Range("M36").Value = 2.85
$ws.Range("N36").Value = 1.93
$ws.Range("O36").Value = 1.7
$ws.Range("P36").Value = 1.39
$ws.Range("Q36").Value = 2.55
$ws.Range("R36").Value = 1.8
$ws.Range("S36").Value = 1.8
$ws.Range("T36").Value = 6.8
$ws.Range("U36").Value = 8.75
$ws.Range("V36").Value = 8.5
$ws.Range("W36").Value = 16.5
$ws.Range("X36").Value = 16
$ws.Range("Y36").Value = 29
$ws.Range("Z36").Value = 9

$ws.Range("AA37").Value = 8.25
$ws.Range("AB37").Value = 19
$ws.Range("AC37").Value = 90
$ws.Range("AD37").Value = 800
$ws.Range("AE37").Value = 6.6
$ws.Range("AF37").Value = 6.6
$ws.Range("AG37").Value = 8.25
$ws.Range("AH37").Value = 9.75
$ws.Range("AI37").Value = 12
$ws.Range("AJ37").Value = 28
$ws.Range("G37").Value = 6.2
$ws.Range("H37").Value = 4.15
$ws.Range("I37").Value = 1.45
$ws.Range("L37").Value = 1.24
$ws.Range("M37").Value = 3.25
$ws.Range("N37").Value = 1.72
$ws.Range("O37").Value = 1.88
$ws.Range("R37").Value = 1.9
$ws.Range("S37").Value = 1.72
$ws.Range("T37").Value = 16
$ws.Range("U37").Value = 40
$ws.Range("V37").Value = 20
$ws.Range("W37").Value = 120
$ws.Range("X37").Value = 70
$ws.Range("Y37").Value = 65
$ws.Range("Z37").Value = 11.5

$ws.Range("AB39").Value = 23
$ws.Range("AE39").Value = 26
$ws.Range("AH39").Value = 151
$ws.Range("AI39").Value = 81
$ws.Range("G39").Value = 1.27
$ws.Range("H39").Value = 5.25
$ws.Range("J39").Value = 1.04
$ws.Range("K39").Value = 13
$ws.Range("W39").Value = 7.5
